$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 19: replace Lysionotus_chingii data with Briggsia_longipes data
$ws.Range("A19").Value = "Briggsia_longipes"
$ws.Range("B19").Value = 3
$ws.Range("C19").Value = 42
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 45

# Row 20 (Briggsia_mihieri) values are unchanged.

# Row 21 (Raphiocarpus_begoniifolius) is removed entirely.
$ws.Rows.Item(21).Delete()

# Update selection to match new active cell
$ws.Range("I20").Select()
